$d = $word.ActiveDocument
$sec = $d.Sections(1)
$h = $sec.Headers(1)
$shp = $h.Range.InlineShapes(1)
$newshp = $shp.ConvertToShape()
Write-Host "Converted. Name:" $newshp.Name
$newshp.Name = "image2.jpg"
Write-Host "Name after set:" $newshp.Name
